$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date value (row 8, column B) to reflect the new generation timestamp.
$ws.Range("B8").Value = "2025-10-02T18:31:12+01:00"

# Set the "Case Sensitive" value (row 20, column B) to "true".
# A leading apostrophe forces Excel to store this as literal text instead of
# auto-converting the word "true" into a native boolean value.
$ws.Range("B20").Value = "'true"

# The apostrophe-quoted entry flips the cell to a "quote prefix" style variant;
# re-copy the plain formatting from the cell above so B20 keeps the same
# unobtrusive style as the rest of the column.
$ws.Range("B19").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
